$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the duplicated "Actual result" (H column) values for rows 3-6,
# keeping their existing cell style/formatting intact.
$ws.Range("H3").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("H6").Value = ""

# Move the active selection to G10 (matches the saved sheet view state).
$ws.Range("G10").Select()
